# Update the cryptos price/volume table with the latest scrape results.
# Numeric-looking "Price" cells (column D) are forced to text (NumberFormat "@")
# before assignment so Excel doesn't reinterpret e.g. "0.620" or "6.23" as numbers
# and silently drop significant trailing zeros / add float noise.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.976.09"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "2.044.15"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.62"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.85"
$ws.Range("E8").Value = "  +3.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.30"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0751"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "2.346.73"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.20"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.67"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.767"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.11"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "2.046.41"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "36.923.68"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.23"
$ws.Range("E20").Value = "  +10.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.60"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").Value = "0.0₃0803"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "224.17"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.35"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.56"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.44"
$ws.Range("E28").Value = "  +6.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.71"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.94"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.44"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.68"
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.42"
$ws.Range("E42").Value = "  +3.32%  "
$ws.Range("D43").Value = "1.476.35"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.68"
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0919"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.15"
$ws.Range("E46").Value = "  +3.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0208"
$ws.Range("E47").Value = "  +2.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.14"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.94"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.08"
$ws.Range("E51").Value = "  +2.52%  "
